$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $old"
    }
}

# --- Title ---
Replace-Exact "The Butterfly Effect: Unraveling Chaos" "Chemistry: Unveiling the Fabric of Matter"

# --- Author name: "Elara Skevington" -> "Dr" + "." + " Morgan Robertson" (3 runs) ---
Replace-Exact "Elara Skevington" "Dr"
# Re-find the "Dr" run we just created (author line, 36pt) and append the rest after it.
$authorPara = $d.Paragraphs(2)
$authorRange = $authorPara.Range
# authorRange currently covers "Dr" + paragraph mark; trim paragraph mark (2 chars: CR)
$drEnd = $authorRange.Start + 2
$insPoint1 = $d.Range($drEnd, $drEnd)
$insPoint1.InsertAfter(".")
$dotEnd = $drEnd + 1
$insPoint2 = $d.Range($dotEnd, $dotEnd)
$insPoint2.InsertAfter(" Morgan Robertson")

# --- Email address ---
Replace-Exact "elara" "morgan"
Replace-Exact "skevington@emailhost" "robertson@hs"
Replace-Exact "com" "edu"

# --- Body paragraph 1 (first block, before first line break) ---
Replace-Exact "In a realm where the interplay of intricate systems shapes our existence, the Butterfly Effect reigns supreme" "In the realm of science, chemistry stands as a profound and encompassing discipline that unveils the intrinsic nature of matter and its captivating transformations"
Replace-Exact " Atmospheric physicist Edward Lorenz coined this captivating concept, seeding the idea that even the slightest disturbance, akin to a butterfly's delicate wingbeats, can catalyze an unpredictable ripple effect with profound implications" " This realm of discovery encompasses the composition, structure, properties, and behavior of substances that constitute our world"
Replace-Exact " This simple analogy encapsulates the complex and interconnected nature of systems, shedding light on the unpredictable dance of chaos" " From the air we breathe to the food we consume, from the medications that heal us to the materials that build our cities, chemistry is an integral thread interwoven into the fabric of our lives"

# --- Body paragraph 1 (second block, after first double line break) ---
Replace-Exact "Engaging in a game of cosmic billiards, subatomic particles weave a tapestry of possibilities" "As we journey into the captivating world of chemistry, we are presented with a tapestry of vibrant colors, intriguing reactions, and fascinating phenomena"
Replace-Exact " As they traverse their quantum landscape, the mere act of observing their trajectories alters their destinies, injecting an element of inherent uncertainty into the core fabric of reality" " This odyssey of exploration takes us on a journey through the atomic and molecular landscapes, where the fundamental building blocks of matter interact in a complex and intricate dance"
Replace-Exact " This fundamental indeterminacy, captured by Heisenberg's Uncertainty Principle, highlights the limitations of our knowledge and our inability to precisely predict these particles' behaviors" " We unravel the mysteries of chemical bonding, uncovering the patterns and forces that govern the formation and behavior of molecules"

# --- Body paragraph 1 (third block, after second double line break) ---
Replace-Exact "Extrapolating this microscopic chaos to a macroscopic realm, we encounter the butterfly's gentle wingbeats generating a cascade of events leading to a hurricane's furious lashing" "With each experiment conducted and each equation solved, we deepen our knowledge of chemistry's profound impact on our world"
Replace-Exact " At the nexus of interconnectedness, seemingly insignificant actions can sow the seeds of profound consequences, evoking awe and intrigue among those who seek to unravel the secrets of our chaotic world" " We gain insight into the processes that shape our planet's ecosystems, unravel the complexities of metabolic pathways that sustain life, and unlock the secrets of materials that drive technological advancements"

# --- Summary header stays "Summary" (unchanged) ---

# --- Summary body ---
Replace-Exact "The Butterfly Effect encapsulates the potent influence of seemingly insignificant actions, the unpredictable nature of complex systems, and the limitations of our knowledge in predicting their outcomes" "In this extensive essay, we have delved into the alluring world of chemistry, unveiling the intricacies of matter's composition and transformation"
Replace-Exact " Rooted in quantum physics and complex systems theory, it manifests in various fields, urging us to embrace uncertainty and acknowledge the inherent unpredictability that underpins our existence" " From the fundamental principles of atomic and molecular interactions to the practical applications in diverse fields, we have explored the vast tapestry of chemistry's influence"

# Now insert the two new runs (". " + new sentence) before the final period run of the Summary paragraph.
$summaryPara = $d.Paragraphs($d.Paragraphs.Count)
$summaryRange = $summaryPara.Range
$summaryText = $summaryRange.Text
# summaryText ends with ".<CR>" -- the very last character before paragraph mark is the final period run's text "."
# Insert point is right before that last "." (2 chars are paragraph mark: CR)
$lastDotStart = $summaryRange.Start + $summaryText.Length - 2
$insPoint3 = $d.Range($lastDotStart, $lastDotStart)
$insPoint3.InsertAfter(".")
$afterDot = $lastDotStart + 1
$insPoint4 = $d.Range($afterDot, $afterDot)
$insPoint4.InsertAfter(" This exhilarating odyssey has instilled in us an appreciation for the elegance and complexity of the natural world, empowering us with a deeper understanding of the world around us")

# --- Add a new empty paragraph at the very end of the body (after the Summary paragraph) ---
$docEnd = $d.Content.End
$endRange = $d.Range($docEnd - 1, $docEnd - 1)
$endRange.InsertParagraphAfter()

Write-Host "Edit complete"
